$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve existing cell formatting/style for the data range, then force
# text format while writing so numeric-looking strings (e.g. "0.999")
# are not auto-converted to numbers by Excel.
$origStyle = $ws.Range("B2:E51").Style
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.154.61"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "2.265.41"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "305.08"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "96.03"
$ws.Range("E6").Value = "  +3.67%  "
$ws.Range("D7").Value = "0.529"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").Value = "35.24"
$ws.Range("E10").Value = "  +8.75%  "
$ws.Range("D11").Value = "0.0793"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").Value = "6.63"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Value = "2.619.17"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "14.38"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").Value = "2.277.51"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "0.792"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").Value = "42.085.20"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("E19").Value = "  -2.15%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "5.98"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").Value = "67.74"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Value = "237.68"
$ws.Range("E23").Value = "  -2.73%  "
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").Value = "1.93"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "23.72"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").Value = "36.92"
$ws.Range("E28").Value = "  +5.59%  "
$ws.Range("D29").Value = "9.51"
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("D31").Value = "159.46"
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").Value = "5.26"
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "3.19"
$ws.Range("E34").Value = "  +5.83%  "
$ws.Range("D35").Value = "0.0738"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").Value = "17.11"
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("D41").Value = "4.05"
$ws.Range("E41").Value = "  +3.54%  "
$ws.Range("E42").Value = "  +8.61%  "
$ws.Range("D43").Value = "1.988.84"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0284"
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "18.84"
$ws.Range("E45").Value = "  -6.23%  "
$ws.Range("E46").Value = "  -4.33%  "
$ws.Range("E47").Value = "  +0.98%  "
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").Value = "72.11"
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("D51").Value = "90.99"
$ws.Range("E51").Value = "  -0.63%  "

# Restore original style/number-format so unaffected formatting is unchanged.
$ws.Range("B2:E51").Style = $origStyle
